# Update "想去人数" (interest count) values in column F on the "展览" sheet
# and on the corresponding rows of the "全部类型" sheet, reflecting the
# output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# Row => New Value, for the "展览" sheet (column F)
$expoUpdates = @{
    7  = 1345
    13 = 157
    18 = 8637
    20 = 6738
    21 = 10875
    27 = 2021
    28 = 202
    29 = 177
    30 = 2203
    34 = 463
}

foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# Row => New Value, for the "全部类型" sheet (column F)
$allUpdates = @{
    10 = 1345
    18 = 157
    23 = 8637
    25 = 6738
    26 = 10875
    36 = 202
    37 = 177
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
